$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.402730226516724
$ws.Range("B1").Value = 1.984633207321167
$ws.Range("C1").Value = 3.018386125564575
$ws.Range("D1").Value = 3.790305852890015
$ws.Range("E1").Value = 0.9802420735359192
